# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / handoff / handback timestamp
# values that get refreshed each time the report is (re)generated.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest HO Xliff Generate Date for first row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 06:20:59"

# --- zh-cn sheet: Correspond Handoff / Handback Datetime for first row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 06:20:48"
$wsZhCn.Range("K2").Value = "2016-09-07 06:21:33"

# --- de-de sheet: Correspond Handoff Datetime for first row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 06:21:50"
